# Generate Report for Handback
# Update the timestamps recorded in the handback status report (the
# "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" columns) to reflect the latest run.
# These values are stored as plain text (formatted like a date, but held
# as a shared string, not a numeric date serial), so assigning a string
# keeps the underlying cell type as text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview.Range("G2").Value = "2016-08-25 03:05:41"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the first file row.
$wsZhCn.Range("H2").Value = "2016-08-25 03:05:36"
$wsZhCn.Range("K2").Value = "2016-08-25 03:05:53"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the first file row.
$wsDeDe.Range("H2").Value = "2016-08-25 03:05:41"
$wsDeDe.Range("K2").Value = "2016-08-25 03:06:02"
